$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# Update the time_taken column (F) with refreshed timestamps
$data.Range("F2").Value = "2021-10-05 14:34:57.016614"
$data.Range("F3").Value = "2021-10-05 14:34:57.016622"
$data.Range("F4").Value = "2021-10-05 14:34:57.016625"
$data.Range("F5").Value = "2021-10-05 14:34:57.016627"
$data.Range("F6").Value = "2021-10-05 14:34:57.016630"
$data.Range("F7").Value = "2021-10-05 14:34:57.016633"
$data.Range("F8").Value = "2021-10-05 14:34:57.016636"

# Add the new metadata sheet after the data sheet
$meta = $wb.Worksheets.Add()
$meta.Name = "metadata"
$meta.Move($null, $data)

# Header row
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"
$meta.Range("B1:G1").Style = $data.Range("B1").Style

# Data row
$meta.Range("A2").Value = 0
$meta.Range("A2").Style = $data.Range("A2").Style
$meta.Range("B2").Value = "Multiple epiphyseal dysplasia and pseudoachondroplasia"
$meta.Range("C2").Value = 3127
$meta.Range("D2").Value = "0.4"
$meta.Range("E2").Value = "2021-09-02T22:34:40.715575Z"
$meta.Range("F2").Value = "2021-10-05 14:34:57.013101"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3127/?format=json"
